# Saldo.xlsx update
# - Row 3 (account 003895497 / EDNA / 100000) becomes account 003301389 /
#   EDMUNDO / 123698.6 (consolidating the old row-32 EDMUNDO balance into
#   this account, with an updated balance).
# - The now-redundant rows are removed (deleted bottom-up so earlier row
#   indices stay valid while we work):
#     row 32: 003301389 EDMUNDO 564.98   (merged up into row 3 above)
#     row 23: 005581299 ZILDA    1458.7
#     row 19: 005009947 VERANICE 5000
#     row 17: 004260002 ERICA    6000
#     row 13: 005133039 PAULO   13168.76
#     row 11: 004268684 PATRICIA 19753.14
#     row 8 : 005440756 VALERIA 46415.02
#     row 7 : 004212409 RAFAEL  49258.71
#     row 6 : 004335144 EDMUNDO 51133.62
# - A new account is added right after the "LEANDRO" row:
#     002687737 / JOSE / 28.13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

# --- 1. Delete the rows that disappear, from the bottom up so the row
#        numbers above the deletion point never shift underneath us. ---
$rowsToDelete = @(32, 23, 19, 17, 13, 11, 8, 7, 6)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# --- 2. Update row 3 in place: 003895497/EDNA/100000 -> 003301389/EDMUNDO/123698.6
Set-TextCell 3 1 "003301389"
Set-TextCell 3 2 "EDMUNDO"
$ws.Cells.Item(3, 3).Value2 = 123698.6

# --- 3. Insert the new JOSE row right after LEANDRO (which, after the
#        deletions above, now sits at row 123), so it lands on row 124. ---
$insertRow = 124
$ws.Rows.Item($insertRow).Insert()
Set-TextCell $insertRow 1 "002687737"
Set-TextCell $insertRow 2 "JOSE"
$ws.Cells.Item($insertRow, 3).Value2 = 28.13
